# "fixed script for run NF"
# The Scene.xml FilePath references were pointing one directory too deep
# (../../NFDataCfg/...) -- fix them to the correct relative path
# (../NFDataCfg/...) used by the NF config loader.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("F12").Value = "../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("F13").Value = "../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("F14").Value = "../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("F15").Value = "../NFDataCfg/Ini/Scene/6.xml"

# Update the active selection to match the saved view state.
$ws.Range("F23").Select()
